$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates - B2 cleared, C2/D2/E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 5.3266291156268153
$ws.Range("D2").Value = 2.3994262606171302
$ws.Range("E2").Value = 3.8708698498221517

# Row 3 updates
$ws.Range("B3").Value = 1.5938355652301459
$ws.Range("C3").Value = 8.0717059160327356
$ws.Range("D3").Value = 6.7690085495513559
$ws.Range("E3").Value = 10.79845978970932

# Update selection to match new active range
$ws.Range("B1:E3").Select()
